$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: station changed from USACE 76065 -> USGS 073802332
$ws.Range("A2").Value = "USGS"
$ws.Range("B2").Value = "073802332"

# Row 7: station 82742 -> 82740, datum offset now -0.16
$ws.Range("B7").Value = "82740"
$ws.Range("C7").Value = -0.16

# Row 11: datum offset now -4.5
$ws.Range("C11").Value = -4.5

# New station row inserted at row 23 (USACE 76305, offset -0.08);
# everything that was row 23-29 shifts down to 24-30.
$ws.Rows(23).Insert()
$ws.Range("A23").Value = "USACE"
$ws.Range("B23").Value = "76305"
$ws.Range("C23").Value = -0.08

# Selection moved to D23 in the saved workbook
$ws.Range("D23").Select()
